$wb = $excel.ActiveWorkbook

# --- Status text changed from "Ready for handoff" to "In Translation" ---
# (Overview!E2/F2 hold the per-locale status for zh-cn/de-de; each locale
#  sheet's Status column, C2, mirrors the same value.)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Status columns narrowed to fit the shorter text ---
$overview.Columns.Item(5).ColumnWidth = 13
$overview.Columns.Item(6).ColumnWidth = 13
$zhcn.Columns.Item(3).ColumnWidth = 13
$dede.Columns.Item(3).ColumnWidth = 13
